# ---------------------------------------------------------------------------
# Weekly CompStat update for the 42nd Precinct report:
#   - Bump "Volume 32 Number 16" -> "Number 17" in the report title.
# - Update the reporting week from 4/14/2025-4/20/2025 to 4/21/2025-4/27/2025.
#   - Refresh the Week-to-Date / 28-Day / Year-to-Date / 2-Year crime stats
#     table (rows 15-30) with newly collected figures.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text runs; only the trailing run's text
#     actually changes, so we target just that substring of each cell) ---

# A8: "Volume 32   Number  16" -> "...  17"
$ws.Cells.Item(8,1).Characters(21,2).Text = "17"

# C9: "Report Covering the Week  4/14/2025  Through  4/20/2025"
#  -> "Report Covering the Week  4/21/2025  Through  4/27/2025"
$ws.Cells.Item(9,3).Characters(27,9).Text = "4/21/2025"
$ws.Cells.Item(9,3).Characters(47,9).Text = "4/27/2025"

# --- Crime-complaints table refresh (rows 15-30) ---

# Row 15
$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 1
$arr[0,1] = 1
$arr[0,2] = 0
$arr[0,3] = 1
$arr[0,4] = 1
$arr[0,5] = 0
$arr[0,6] = 11
$arr[0,7] = 15
$arr[0,8] = -26.666666666666
$arr[0,9] = -26.666666666666
$arr[0,10] = 120
$arr[0,11] = -45
$ws.Range("C15:N15").Value = $arr
$ws.Range("C15:D15").NumberFormat = $ws.Range("C16:D16").NumberFormat
$ws.Range("E15").NumberFormat = $ws.Range("E16").NumberFormat

# Row 16
$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 12
$arr[0,1] = 12
$arr[0,2] = 0
$arr[0,3] = 39
$arr[0,4] = 27
$arr[0,5] = 44.444444444444
$arr[0,6] = 136
$arr[0,7] = 128
$arr[0,8] = 6.25
$arr[0,9] = 12.396694214876
$arr[0,10] = 76.623376623376
$arr[0,11] = -63.538873994638
$ws.Range("C16:N16").Value = $arr

# Row 17
$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 9
$arr[0,1] = 15
$arr[0,2] = -40
$arr[0,3] = 45
$arr[0,4] = 62
$arr[0,5] = -27.419354838709
$arr[0,6] = 200
$arr[0,7] = 224
$arr[0,8] = -10.714285714285
$arr[0,9] = 5.820105820105
$arr[0,10] = 156.410256410256
$arr[0,11] = -28.057553956834
$ws.Range("C17:N17").Value = $arr

# Row 18
$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 3
$arr[0,1] = 9
$arr[0,2] = -66.666666666666
$arr[0,3] = 24
$arr[0,4] = 23
$arr[0,5] = 4.347826086956
$arr[0,6] = 81
$arr[0,7] = 93
$arr[0,8] = -12.903225806451
$arr[0,9] = 19.117647058823
$arr[0,10] = 113.157894736842
$arr[0,11] = -75.964391691394
$ws.Range("C18:N18").Value = $arr

# Row 19
$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 11
$arr[0,1] = 10
$arr[0,2] = 10
$arr[0,3] = 40
$arr[0,4] = 44
$arr[0,5] = -9.090909090909
$arr[0,6] = 162
$arr[0,7] = 189
$arr[0,8] = -14.285714285714
$arr[0,9] = 14.084507042253
$arr[0,10] = 131.428571428571
$arr[0,11] = 70.526315789473
$ws.Range("C19:N19").Value = $arr

# Row 20
$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 5
$arr[0,1] = 6
$arr[0,2] = -16.666666666666
$arr[0,3] = 26
$arr[0,4] = 16
$arr[0,5] = 62.5
$arr[0,6] = 84
$arr[0,7] = 76
$arr[0,8] = 10.526315789473
$arr[0,9] = -44
$arr[0,10] = 133.333333333333
$arr[0,11] = -55.080213903743
$ws.Range("C20:N20").Value = $arr

# Row 21
$arr = New-Object 'object[,]' 1,12
$arr[0,0] = 41
$arr[0,1] = 53
$arr[0,2] = -22.641509433962
$arr[0,3] = 175
$arr[0,4] = 173
$arr[0,5] = 1.156069364161
$arr[0,6] = 674
$arr[0,7] = 726
$arr[0,8] = -7.162534435261
$arr[0,9] = -2.177068214804
$arr[0,10] = 118.122977346278
$arr[0,11] = -48.073959938366
$ws.Range("C21:N21").Value = $arr

# Row 22
$ws.Range("M22").Value = -71.428571428571

# Row 23
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 4
$arr[0,1] = 10
$arr[0,2] = -60
$arr[0,3] = 25
$arr[0,4] = 26
$arr[0,5] = -3.846153846153
$arr[0,6] = 86
$arr[0,7] = 129
$arr[0,8] = -33.333333333333
$arr[0,9] = -39.007092198581
$arr[0,10] = 48.275862068965
$ws.Range("C23:M23").Value = $arr

# Row 24
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 31
$arr[0,1] = 27
$arr[0,2] = 14.814814814814
$arr[0,3] = 110
$arr[0,4] = 91
$arr[0,5] = 20.87912087912
$arr[0,6] = 384
$arr[0,7] = 364
$arr[0,8] = 5.494505494505
$arr[0,9] = 6.077348066298
$arr[0,10] = 67.68558951965
$ws.Range("C24:M24").Value = $arr

# Row 25
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 12
$arr[0,1] = 2
$arr[0,2] = 500
$arr[0,3] = 29
$arr[0,4] = 13
$arr[0,5] = 123.076923076923
$arr[0,6] = 77
$arr[0,7] = 54
$arr[0,8] = 42.592592592592
$arr[0,9] = 1.315789473684
$ws.Range("C25:L25").Value = $arr

# Row 26
$arr = New-Object 'object[,]' 1,11
$arr[0,0] = 18
$arr[0,1] = 19
$arr[0,2] = -5.263157894736
$arr[0,3] = 74
$arr[0,4] = 69
$arr[0,5] = 7.246376811594
$arr[0,6] = 278
$arr[0,7] = 273
$arr[0,8] = 1.831501831501
$arr[0,9] = -24.043715846994
$arr[0,10] = 3.731343283582
$ws.Range("C26:M26").Value = $arr

# Row 27
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1
$arr[0,1] = 1
$arr[0,2] = 0
$arr[0,3] = 1
$arr[0,4] = 3
$arr[0,5] = -66.666666666666
$arr[0,6] = 12
$arr[0,7] = 21
$arr[0,8] = -42.857142857142
$arr[0,9] = -45.454545454545
$ws.Range("C27:L27").Value = $arr
$ws.Range("C27:D27").NumberFormat = $ws.Range("C26:D26").NumberFormat
$ws.Range("E27").NumberFormat = $ws.Range("E26").NumberFormat

# Row 28
$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1
$arr[0,1] = 2
$arr[0,2] = -50
$arr[0,3] = 7
$arr[0,4] = 6
$arr[0,5] = 16.666666666666
$arr[0,6] = 17
$arr[0,7] = 23
$arr[0,8] = -26.086956521739
$arr[0,9] = -48.484848484848
$ws.Range("C28:L28").Value = $arr

# Row 29
$ws.Range("D29").Value = 2
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 3
$arr[0,1] = -66.666666666666
$ws.Range("G29:H29").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 6
$arr[0,1] = 16.666666666666
$ws.Range("J29:K29").Value = $arr
$ws.Range("N29").Value = -70.833333333333

# Row 30
$ws.Range("D30").Value = 2
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 3
$arr[0,1] = -66.666666666666
$ws.Range("G30:H30").Value = $arr
$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 6
$arr[0,1] = 0
$ws.Range("J30:K30").Value = $arr
$ws.Range("N30").Value = -75

